$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 40962.56
$ws.Range("I28").Value = 48241.668
$ws.Range("J28").Value = 2747.25
$ws.Range("K28").Value = 48241.668
$ws.Range("L28").Value = 2747.25
$ws.Range("M28").Value = -47756.668
$ws.Range("N28").Value = -3717.25

$ws.Range("H62").Value = 755984.9399999999
$ws.Range("I62").Value = 2004400.8
$ws.Range("J62").Value = 62420.555
$ws.Range("K62").Value = 2004400.8
$ws.Range("L62").Value = 62420.555
$ws.Range("M62").Value = -2003776.8
$ws.Range("N62").Value = -63668.555

$ws.Range("H64").Value = 4088.0908
$ws.Range("I64").Value = 3999.5
$ws.Range("J64").Value = 4107.778
$ws.Range("K64").Value = 3999.5
$ws.Range("L64").Value = 4107.778
$ws.Range("M64").Value = -3751.5
$ws.Range("N64").Value = -4603.778

$ws.Range("H65").Value = 755984.9399999999
$ws.Range("I65").Value = 2004400.8
$ws.Range("J65").Value = 62420.555
$ws.Range("K65").Value = 10022004
$ws.Range("L65").Value = 312102.775
$ws.Range("M65").Value = -10018884
$ws.Range("N65").Value = -318342.775

$ws.Range("H67").Value = 4088.0908
$ws.Range("I67").Value = 3999.5
$ws.Range("J67").Value = 4107.778
$ws.Range("K67").Value = 3999.5
$ws.Range("L67").Value = 4107.778
$ws.Range("M67").Value = -3141.5
$ws.Range("N67").Value = -5823.778

$ws.Range("H74").Value = 3800
$ws.Range("I74").Value = 3133.3333
$ws.Range("K74").Value = 3133.3333
$ws.Range("M74").Value = -2197.3333

$ws.Range("H77").Value = 3800
$ws.Range("I77").Value = 3133.3333
$ws.Range("K77").Value = 15666.6665
$ws.Range("M77").Value = -10986.6665

$ws.Range("H86").Value = 7166862
$ws.Range("I86").Value = 4775.4
$ws.Range("J86").Value = 15430808
$ws.Range("K86").Value = 4775.4
$ws.Range("L86").Value = 15430808
$ws.Range("M86").Value = -3652.4
$ws.Range("N86").Value = -15433054

$ws.Range("H89").Value = 7166862
$ws.Range("I89").Value = 4775.4
$ws.Range("J89").Value = 15430808
$ws.Range("K89").Value = 23877
$ws.Range("L89").Value = 77154040
$ws.Range("M89").Value = -18261
$ws.Range("N89").Value = -77165272

$ws.Range("H97").Value = 1083.3334
$ws.Range("J97").Value = 1325
$ws.Range("L97").Value = 3975
$ws.Range("N97").Value = -4967

$ws.Range("H106").Value = 10863.077
$ws.Range("I106").Value = 6308.875
$ws.Range("K106").Value = 6308.875
$ws.Range("M106").Value = -5677.875

$ws.Range("H107").Value = 696.05884
$ws.Range("I107").Value = 697.5333000000001
$ws.Range("J107").Value = 685
$ws.Range("K107").Value = 697.5333000000001
$ws.Range("L107").Value = 685
$ws.Range("M107").Value = 1222.4667
$ws.Range("N107").Value = -4525

$ws.Range("H112").Value = 2088
$ws.Range("I112").Value = 2700
$ws.Range("J112").Value = 1959.1578
$ws.Range("K112").Value = 8100
$ws.Range("L112").Value = 5877.4734
$ws.Range("M112").Value = -6992
$ws.Range("N112").Value = -8093.4734

$ws.Range("H115").Value = 493
$ws.Range("I115").Value = 408.5
$ws.Range("K115").Value = 1225.5
$ws.Range("M115").Value = 341.5

$ws.Range("H116").Value = 18267446
$ws.Range("I116").Value = 12555435
$ws.Range("J116").Value = 27787462
$ws.Range("K116").Value = 12555435
$ws.Range("L116").Value = 27787462
$ws.Range("M116").Value = -12551993
$ws.Range("N116").Value = -27794346

$ws.Range("H121").Value = 2073.5454
$ws.Range("J121").Value = 2093.4443
$ws.Range("L121").Value = 6280.3329
$ws.Range("N121").Value = -9774.332900000001

$ws.Range("H135").Value = 35715788
$ws.Range("I135").Value = 40000824
$ws.Range("K135").Value = 360007416
$ws.Range("M135").Value = -360004881

$ws.Range("H137").Value = 3860.02
$ws.Range("I137").Value = 2793.5
$ws.Range("J137").Value = 3904.4583
$ws.Range("K137").Value = 8380.5
$ws.Range("L137").Value = 11713.3749
$ws.Range("M137").Value = -5830.5
$ws.Range("N137").Value = -16813.3749

$ws.Range("H138").Value = 3076.4285
$ws.Range("I138").Value = 2460.375
$ws.Range("J138").Value = 3897.8333
$ws.Range("K138").Value = 7381.125
$ws.Range("L138").Value = 11693.4999
$ws.Range("M138").Value = -2241.125
$ws.Range("N138").Value = -21973.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 175929.14
$ws.Range("I32").Value = 190505.67
$ws.Range("J32").Value = 47169.668
$ws.Range("K32").Value = 190505.67
$ws.Range("L32").Value = 47169.668
$ws.Range("M32").Value = -190218.67
$ws.Range("N32").Value = -47743.668

$ws.Range("H61").Value = 3159.4
$ws.Range("I61").Value = 3159.4
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3159.4
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2947.4
$ws.Range("N61").Value = $null

$ws.Range("H74").Value = 2287.4348
$ws.Range("I74").Value = 2130.6
$ws.Range("K74").Value = 2130.6
$ws.Range("M74").Value = -1256.6

$ws.Range("H77").Value = 2287.4348
$ws.Range("I77").Value = 2130.6
$ws.Range("K77").Value = 10653
$ws.Range("M77").Value = -6285

$ws.Range("H136").Value = 3159.4
$ws.Range("I136").Value = 3159.4
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9478.200000000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6928.200000000001
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5000
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5226

$ws.Range("H55").Value = 70178.664
$ws.Range("I55").Value = 69558
$ws.Range("J55").Value = 71420
$ws.Range("K55").Value = 69558
$ws.Range("L55").Value = 71420
$ws.Range("M55").Value = -69285
$ws.Range("N55").Value = -71966

$ws.Range("H86").Value = 25002158
$ws.Range("I86").Value = 38463388
$ws.Range("J86").Value = 2726.8572
$ws.Range("K86").Value = 38463388
$ws.Range("L86").Value = 2726.8572
$ws.Range("M86").Value = -38462265
$ws.Range("N86").Value = -4972.8572

$ws.Range("H89").Value = 25002158
$ws.Range("I89").Value = 38463388
$ws.Range("J89").Value = 2726.8572
$ws.Range("K89").Value = 192316940
$ws.Range("L89").Value = 13634.286
$ws.Range("M89").Value = -192311324
$ws.Range("N89").Value = -24866.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3755.4321
$ws.Range("I31").Value = 1097.0769
$ws.Range("J31").Value = 5012.109
$ws.Range("K31").Value = 1097.0769
$ws.Range("L31").Value = 5012.109
$ws.Range("M31").Value = -802.0769
$ws.Range("N31").Value = -5602.109

$ws.Range("H34").Value = 3755.4321
$ws.Range("I34").Value = 1097.0769
$ws.Range("J34").Value = 5012.109
$ws.Range("K34").Value = 1097.0769
$ws.Range("L34").Value = 5012.109
$ws.Range("M34").Value = -895.0769
$ws.Range("N34").Value = -5416.109

$ws.Range("H68").Value = 69991
$ws.Range("J68").Value = 69991
$ws.Range("L68").Value = 69991
$ws.Range("N68").Value = -71489

$ws.Range("H71").Value = 69991
$ws.Range("J71").Value = 69991
$ws.Range("L71").Value = 209973
$ws.Range("N71").Value = -217461

$ws.Range("H134").Value = 3061.6
$ws.Range("I134").Value = 2101.9
$ws.Range("K134").Value = 6305.700000000001
$ws.Range("M134").Value = -3770.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 49001
$ws.Range("J39").Value = 49001
$ws.Range("L39").Value = 49001
$ws.Range("N39").Value = -50065

$ws.Range("H43").Value = 4798.8
$ws.Range("I43").Value = 4798.8
$ws.Range("K43").Value = 4798.8
$ws.Range("M43").Value = -4647.8

$ws.Range("H46").Value = 35999.75
$ws.Range("J46").Value = 49499.5
$ws.Range("L46").Value = 49499.5
$ws.Range("N46").Value = -49811.5

$ws.Range("H80").Value = 2870
$ws.Range("I80").Value = 2572.3635
$ws.Range("J80").Value = 3961.3333
$ws.Range("K80").Value = 2572.3635
$ws.Range("L80").Value = 3961.3333
$ws.Range("M80").Value = -1574.3635
$ws.Range("N80").Value = -5957.3333

$ws.Range("H83").Value = 2870
$ws.Range("I83").Value = 2572.3635
$ws.Range("J83").Value = 3961.3333
$ws.Range("K83").Value = 12861.8175
$ws.Range("L83").Value = 19806.6665
$ws.Range("M83").Value = -7869.817499999999
$ws.Range("N83").Value = -29790.6665

$ws.Range("H113").Value = 3162.95
$ws.Range("J113").Value = 3898.6667
$ws.Range("L113").Value = 3898.6667
$ws.Range("N113").Value = -8238.6667

$ws.Range("H132").Value = 235893.16
$ws.Range("I132").Value = 306016.8
$ws.Range("J132").Value = 4485.1
$ws.Range("K132").Value = 918050.3999999999
$ws.Range("L132").Value = 13455.3
$ws.Range("M132").Value = -915520.3999999999
$ws.Range("N132").Value = -18515.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 13249.5
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251

$ws.Range("H71").Value = 13249.5
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

$ws.Range("H136").Value = 2529.8
$ws.Range("I136").Value = 1787.375
$ws.Range("J136").Value = 5499.5
$ws.Range("K136").Value = 5362.125
$ws.Range("L136").Value = 16498.5
$ws.Range("M136").Value = -2812.125
$ws.Range("N136").Value = -21598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2389340.5
$ws.Range("I62").Value = 7949137
$ws.Range("J62").Value = 6570.7144
$ws.Range("K62").Value = 7949137
$ws.Range("L62").Value = 6570.7144
$ws.Range("M62").Value = -7948513
$ws.Range("N62").Value = -7818.7144

$ws.Range("H65").Value = 2389340.5
$ws.Range("I65").Value = 7949137
$ws.Range("J65").Value = 6570.7144
$ws.Range("K65").Value = 39745685
$ws.Range("L65").Value = 32853.572
$ws.Range("M65").Value = -39742565
$ws.Range("N65").Value = -39093.572

$ws.Range("H122").Value = 1557
$ws.Range("I122").Value = 1591.04
$ws.Range("K122").Value = 4773.12
$ws.Range("M122").Value = -2323.12
